# Weekly update: insert a new price record as the new row 9 (pushing all
# existing records for this product down by one row), matching the
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 9 (existing rows 9..57 shift
# down to 10..58; column D's date style is carried onto the new row
# automatically, same as native Excel row-insert behaviour).
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44677
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107011
$ws.Range("J9").Value = "Tuna"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 55
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = "$/caja 16 kilos"
$ws.Range("R9").Value = "Provincia de Los Andes"
$ws.Range("S9").Value = 1250
$ws.Range("T9").Value = 16
